$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = '@'
$cell.Value = '28.672.53'
$cell.Style = 'Normal'
$ws.Range("E2").Value = '  +1.17%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = '@'
$cell.Value = '1.566.57'
$cell.Style = 'Normal'
$ws.Range("E3").Value = '  -0.05%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range("E4").Value = '  -0.25%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = '@'
$cell.Value = '210.09'
$cell.Style = 'Normal'
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("E6").Value = '  -0.52%  '

$ws.Range("E7").Value = '  -0.28%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = '@'
$cell.Value = '25.06'
$cell.Style = 'Normal'
$ws.Range("E8").Value = '  +5.56%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = '@'
$cell.Value = '0.244'
$cell.Style = 'Normal'
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("E10").Value = '  -0.07%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = '@'
$cell.Value = '0.0896'
$cell.Style = 'Normal'
$ws.Range("E11").Value = '  +0.36%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = '@'
$cell.Value = '1.792.06'
$cell.Style = 'Normal'
$ws.Range("E12").Value = '  +0.10%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = '@'
$cell.Value = '1.562.93'
$cell.Style = 'Normal'
$ws.Range("E13").Value = '  -0.30%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = '@'
$cell.Value = '28.694.06'
$cell.Style = 'Normal'
$ws.Range("E14").Value = '  +1.31%  '

$ws.Range("E15").Value = '  +0.84%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = '@'
$cell.Value = '3.64'
$cell.Style = 'Normal'
$ws.Range("E16").Value = '  -0.60%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = '@'
$cell.Value = '61.38'
$cell.Style = 'Normal'
$ws.Range("E17").Value = '  +0.47%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = '@'
$cell.Value = '229.44'
$cell.Style = 'Normal'
$ws.Range("E18").Value = '  +0.76%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = '@'
$cell.Value = '7.35'
$cell.Style = 'Normal'
$ws.Range("E19").Value = '  -0.48%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0679'
$cell.Style = 'Normal'
$ws.Range("E20").Value = '  +0.09%  '

$cell = $ws.Range("D21")
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'
$ws.Range("E21").Value = '  -0.30%  '

$ws.Range("E22").Value = '  -0.25%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = '@'
$cell.Value = '9.02'
$cell.Style = 'Normal'
$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("E24").Value = '  +1.05%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = '@'
$cell.Value = '151.92'
$cell.Style = 'Normal'
$ws.Range("E25").Value = '  +0.94%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = '@'
$cell.Value = '14.79'
$cell.Style = 'Normal'
$ws.Range("E26").Value = '  -0.66%  '

$ws.Range("E27").Value = '  +0.30%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range("E28").Value = '  -0.25%  '

$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("E30").Value = '  -4.02%  '

$ws.Range("E31").Value = '  -2.56%  '

$ws.Range("E32").Value = '  +0.16%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = '@'
$cell.Value = '1.394.39'
$cell.Style = 'Normal'
$ws.Range("E33").Value = '  +0.95%  '

$ws.Range("E34").Value = '  -2.84%  '

$ws.Range("E35").Value = '  -3.87%  '

$ws.Range("E36").Value = '  -1.36%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = '@'
$cell.Value = '2.70'
$cell.Style = 'Normal'
$ws.Range("E37").Value = '  +1.84%  '

$ws.Range("E38").Value = '  -2.14%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = '@'
$cell.Value = '0.0161'
$cell.Style = 'Normal'
$ws.Range("E39").Value = '  -0.59%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = '@'
$cell.Value = '1.97'
$cell.Style = 'Normal'
$ws.Range("E40").Value = '  +1.95%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = '@'
$cell.Value = '0.521'
$cell.Style = 'Normal'
$ws.Range("E41").Value = '  -0.06%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range("E42").Value = '  -0.24%  '

$ws.Range("E43").Value = '  -1.41%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = '@'
$cell.Value = '0.0460'
$cell.Style = 'Normal'
$ws.Range("E44").Value = '  -3.46%  '

$ws.Range("E45").Value = '  +2.78%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = '@'
$cell.Value = '5.24'
$cell.Style = 'Normal'
$ws.Range("E46").Value = '  -1.50%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = '@'
$cell.Value = '1.702.72'
$cell.Style = 'Normal'
$ws.Range("E47").Value = '  +0.04%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = '@'
$cell.Value = '0.870'
$cell.Style = 'Normal'
$ws.Range("E48").Value = '  -5.11%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = '@'
$cell.Value = '85.09'
$cell.Style = 'Normal'
$ws.Range("E49").Value = '  -0.57%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = '@'
$cell.Value = '43.00'
$cell.Style = 'Normal'
$ws.Range("E50").Value = '  +5.46%  '

$ws.Range("E51").Value = '  -0.11%  '
